$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.182.41"
$ws.Range("E2").Value = "  -0.50%  "
$ws.Range("D3").Value = "3.506.65"
$ws.Range("E3").Value = "  -0.93%  "
$ws.Range("D5").Value = "606.15"
$ws.Range("E5").Value = "  -0.24%  "
$ws.Range("D6").Value = "172.86"
$ws.Range("E6").Value = "  -0.57%  "
$ws.Range("E7").Value = "  -1.42%  "
$ws.Range("D8").Value = "3.502.70"
$ws.Range("E8").Value = "  -0.94%  "
$ws.Range("E9").Value = "  -0.05%  "
$ws.Range("E10").Value = "  -2.29%  "
$ws.Range("D11").Value = "7.28"
$ws.Range("E11").Value = "  +7.65%  "
$ws.Range("D12").Value = "0.586"
$ws.Range("E12").Value = "  +0.03%  "
$ws.Range("D13").Value = "46.15"
$ws.Range("E13").Value = "  -2.89%  "
$ws.Range("E14").Value = "  -1.35%  "
$ws.Range("D15").Value = "4.077.12"
$ws.Range("E15").Value = "  -0.83%  "
$ws.Range("D16").Value = "8.35"
$ws.Range("E16").Value = "  -1.31%  "
$ws.Range("D17").Value = "614.47"
$ws.Range("E17").Value = "  -2.24%  "
$ws.Range("D18").Value = "3.504.21"
$ws.Range("E18").Value = "  -1.23%  "
$ws.Range("D19").Value = "70.194.01"
$ws.Range("E19").Value = "  -0.59%  "
$ws.Range("E20").Value = "  +1.09%  "
$ws.Range("D21").Value = "17.48"
$ws.Range("E21").Value = "  +0.25%  "
$ws.Range("D22").Value = "0.879"
$ws.Range("E22").Value = "  -1.08%  "
$ws.Range("D23").Value = "9.15"
$ws.Range("E23").Value = "  -8.62%  "
$ws.Range("D24").Value = "98.42"
$ws.Range("E24").Value = "  +1.45%  "
$ws.Range("D25").Value = "15.57"
$ws.Range("E25").Value = "  -1.84%  "
$ws.Range("E26").Value = "  -3.49%  "
$ws.Range("E27").Value = "  -0.05%  "
$ws.Range("E28").Value = "  -1.90%  "
$ws.Range("D29").Value = "33.84"
$ws.Range("E29").Value = "  +1.37%  "
$ws.Range("D30").Value = "8.99"
$ws.Range("E30").Value = "  -2.18%  "
$ws.Range("D31").Value = "2.99"
$ws.Range("E31").Value = "  -3.99%  "
$ws.Range("D32").Value = "8.05"
$ws.Range("E32").Value = "  -4.86%  "
$ws.Range("E33").Value = "  -4.64%  "
$ws.Range("B34").Value = "Bittensor"
$ws.Range("C34").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D34").Value = "633.06"
$ws.Range("E34").Value = "  +11.91%  "
$ws.Range("B35").Value = "NEARProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D35").Value = "6.81"
$ws.Range("E35").Value = "  -2.99%  "
$ws.Range("D36").Value = "0.0997"
$ws.Range("E36").Value = "  -2.34%  "
$ws.Range("E37").Value = "  -0.19%  "
$ws.Range("D38").Value = "0.0483"
$ws.Range("E38").Value = "  +6.30%  "
$ws.Range("E39").Value = "  -4.64%  "
$ws.Range("D40").Value = "56.75"
$ws.Range("E40").Value = "  -1.19%  "
$ws.Range("E41").Value = "  -0.09%  "
$ws.Range("E42").Value = "  +0.59%  "
$ws.Range("D43").Value = "3.368.77"
$ws.Range("E43").Value = "  +1.20%  "
$ws.Range("E44").Value = "  +2.67%  "
$ws.Range("E45").Value = "  -5.55%  "
$ws.Range("E46").Value = "  -5.14%  "
$ws.Range("D47").Value = "31.94"
$ws.Range("E47").Value = "  -3.31%  "
$ws.Range("E48").Value = "  -3.73%  "
$ws.Range("E49").Value = "  +0.42%  "
$ws.Range("D50").Value = "133.15"
$ws.Range("E50").Value = "  -0.80%  "
